$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1646
$ws.Range("I40").Value = 1562.375
$ws.Range("K40").Value = 1562.375
$ws.Range("M40").Value = -1387.375
$ws.Range("H53").Value = 242.5625
$ws.Range("I53").Value = 175.66667
$ws.Range("K53").Value = 175.66667
$ws.Range("M53").Value = 461.33333
$ws.Range("H69").Value = 29052.79
$ws.Range("J69").Value = 19374.562
$ws.Range("L69").Value = 58123.686
$ws.Range("N69").Value = -59871.686
$ws.Range("H72").Value = 29052.79
$ws.Range("J72").Value = 19374.562
$ws.Range("L72").Value = 174371.058
$ws.Range("N72").Value = -183107.058
$ws.Range("H80").Value = 3092
$ws.Range("I80").Value = 3514.5715
$ws.Range("J80").Value = 2823.0908
$ws.Range("K80").Value = 10543.7145
$ws.Range("L80").Value = 8469.2724
$ws.Range("M80").Value = -9545.7145
$ws.Range("N80").Value = -10465.2724
$ws.Range("H83").Value = 3092
$ws.Range("I83").Value = 3514.5715
$ws.Range("J83").Value = 2823.0908
$ws.Range("K83").Value = 31631.1435
$ws.Range("L83").Value = 25407.8172
$ws.Range("M83").Value = -26639.1435
$ws.Range("N83").Value = -35391.8172
$ws.Range("H106").Value = 4662.5
$ws.Range("I106").Value = 5085.7144
$ws.Range("K106").Value = 5085.7144
$ws.Range("M106").Value = -4454.7144
$ws.Range("H116").Value = 8900
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1675896.1
$ws.Range("I32").Value = 1893858.5
$ws.Range("K32").Value = 1893858.5
$ws.Range("M32").Value = -1893571.5
$ws.Range("H74").Value = 826
$ws.Range("I74").Value = 826
$ws.Range("K74").Value = 826
$ws.Range("M74").Value = 48
$ws.Range("H77").Value = 826
$ws.Range("I77").Value = 826
$ws.Range("K77").Value = 4130
$ws.Range("M77").Value = 238
$ws.Range("H88").Value = 3611.5
$ws.Range("J88").Value = 4199
$ws.Range("L88").Value = 4199
$ws.Range("N88").Value = -5011
$ws.Range("H91").Value = 3611.5
$ws.Range("J91").Value = 4199
$ws.Range("L91").Value = 4199
$ws.Range("N91").Value = -7007
$ws.Range("H97").Value = 1691.5
$ws.Range("J97").Value = 1897.6666
$ws.Range("L97").Value = 1897.6666
$ws.Range("N97").Value = -2889.6666

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 600
$ws.Range("J11").Value = 600
$ws.Range("L11").Value = 600
$ws.Range("N11").Value = -880
$ws.Range("H133").Value = 74999
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 46322.75
$ws.Range("I133").Value = 32648
$ws.Range("K133").Value = 32648
$ws.Range("M133").Value = -30118

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 618.5
$ws.Range("I132").Value = 578.25
$ws.Range("J132").Value = 699
$ws.Range("K132").Value = 5204.25
$ws.Range("L132").Value = 6291
$ws.Range("M132").Value = -2674.25
$ws.Range("N132").Value = -11351

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470
$ws.Range("H140").Value = 142831
$ws.Range("J140").Value = 142831
$ws.Range("L140").Value = 142831
$ws.Range("N140").Value = -153191

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8346.200000000001
$ws.Range("I7").Value = 3234.6667
$ws.Range("J7").Value = 9624.083000000001
$ws.Range("K7").Value = 3234.6667
$ws.Range("L7").Value = 9624.083000000001
$ws.Range("M7").Value = -3122.6667
$ws.Range("N7").Value = -9848.083000000001
$ws.Range("H9").Value = 435
$ws.Range("I9").Value = 55
$ws.Range("J9").Value = 1005
$ws.Range("K9").Value = 55
$ws.Range("L9").Value = 1005
$ws.Range("M9").Value = 169
$ws.Range("N9").Value = -1453
$ws.Range("H22").Value = 891.1667
$ws.Range("I22").Value = 999.25
$ws.Range("J22").Value = 675
$ws.Range("K22").Value = 999.25
$ws.Range("L22").Value = 675
$ws.Range("M22").Value = -704.25
$ws.Range("N22").Value = -1265
$ws.Range("H27").Value = 891.1667
$ws.Range("I27").Value = 999.25
$ws.Range("J27").Value = 675
$ws.Range("K27").Value = 999.25
$ws.Range("L27").Value = 675
$ws.Range("M27").Value = -892.25
$ws.Range("N27").Value = -889
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H61").Value = 2926.1428
$ws.Range("I61").Value = 2926.1428
$ws.Range("K61").Value = 2926.1428
$ws.Range("M61").Value = -2724.1428
$ws.Range("H113").Value = 2926.1428
$ws.Range("I113").Value = 2926.1428
$ws.Range("K113").Value = 2926.1428
$ws.Range("M113").Value = -756.1428000000001
$ws.Range("H126").Value = 8346.200000000001
$ws.Range("I126").Value = 3234.6667
$ws.Range("J126").Value = 9624.083000000001
$ws.Range("K126").Value = 9704.000100000001
$ws.Range("L126").Value = 28872.249
$ws.Range("M126").Value = -7234.000100000001
$ws.Range("N126").Value = -33812.249
$ws.Range("H132").Value = 5162.1
$ws.Range("I132").Value = 5326.6875
$ws.Range("K132").Value = 15980.0625
$ws.Range("M132").Value = -13450.0625

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 1996.3334
$ws.Range("I12").Value = 1996.3334
$ws.Range("K12").Value = 1996.3334
$ws.Range("M12").Value = -1854.3334
$ws.Range("H41").Value = 19982.4
$ws.Range("I41").Value = 19975.5
$ws.Range("J41").Value = 19987
$ws.Range("K41").Value = 19975.5
$ws.Range("L41").Value = 19987
$ws.Range("M41").Value = -19585.5
$ws.Range("N41").Value = -20767
$ws.Range("H74").Value = 29499.5
$ws.Range("J74").Value = 29499.5
$ws.Range("L74").Value = 29499.5
$ws.Range("N74").Value = -31371.5
$ws.Range("H77").Value = 29499.5
$ws.Range("J77").Value = 29499.5
$ws.Range("L77").Value = 88498.5
$ws.Range("N77").Value = -97858.5
$ws.Range("H141").Value = 182149.17
$ws.Range("J141").Value = 175779.2
$ws.Range("L141").Value = 175779.2
$ws.Range("N141").Value = -186139.2
